$wb = $excel.ActiveWorkbook

# --- Sheet: Octubre ---
$ws1 = $wb.Worksheets.Item("Octubre")
$ws1.Range("B2").Value = 2163
$ws1.Range("C2").Value = 982
$ws1.Range("D2").Value = 309
$ws1.Range("E2").Value = 872
$ws1.Range("B3").Value = 11586
$ws1.Range("C3").Value = 610
$ws1.Range("D3").Value = 2381
$ws1.Range("E3").Value = 3001
$ws1.Range("F3").Value = 5594
$ws1.Range("B4").Value = 1448
$ws1.Range("C4").Value = 605
$ws1.Range("D4").Value = 392
$ws1.Range("E4").Value = 451
$ws1.Range("B5").Value = 18543
$ws1.Range("C5").Value = 231
$ws1.Range("D5").Value = 1086
$ws1.Range("E5").Value = 476
$ws1.Range("F5").Value = 16750
$ws1.Range("G5").Value = 17.64
$ws1.Range("H5").Value = 2.22
$ws1.Range("I5").Value = 13
$ws1.Range("B6").Value = 9872
$ws1.Range("C6").Value = 1044
$ws1.Range("D6").Value = 262
$ws1.Range("E6").Value = 878
$ws1.Range("F6").Value = 7688
$ws1.Range("I6").Value = 9
$ws1.Range("J6").Value = 39
$ws1.Range("B7").Value = 17567
$ws1.Range("C7").Value = 1747
$ws1.Range("D7").Value = 327
$ws1.Range("E7").Value = 1931
$ws1.Range("F7").Value = 13562
$ws1.Range("G7").Value = 23.75
$ws1.Range("H7").Value = 2.72
$ws1.Range("I7").Value = 32
$ws1.Range("B8").Value = 20325
$ws1.Range("C8").Value = 1307
$ws1.Range("D8").Value = 1086
$ws1.Range("E8").Value = 4313
$ws1.Range("F8").Value = 13619
$ws1.Range("B9").Value = 13828
$ws1.Range("C9").Value = 3078
$ws1.Range("E9").Value = 479
$ws1.Range("F9").Value = 10202
$ws1.Range("B10").Value = 15629
$ws1.Range("C10").Value = 876
$ws1.Range("D10").Value = 308
$ws1.Range("E10").Value = 406
$ws1.Range("F10").Value = 14039
$ws1.Range("I10").Value = 2
$ws1.Range("B11").Value = 13113
$ws1.Range("C11").Value = 1219
$ws1.Range("D11").Value = 601
$ws1.Range("E11").Value = 2231
$ws1.Range("F11").Value = 9060
$ws1.Range("G11").Value = 12.11
$ws1.Range("H11").Value = 8.199999999999999
$ws1.Range("I11").Value = 51
$ws1.Range("B12").Value = 14681
$ws1.Range("C12").Value = 717
$ws1.Range("D12").Value = 552
$ws1.Range("E12").Value = 586
$ws1.Range("F12").Value = 12825
$ws1.Range("G12").Value = 9.19
$ws1.Range("H12").Value = 13.36
$ws1.Range("I12").Value = 24
$ws1.Range("J12").Value = 39
$ws1.Range("B13").Value = 7995
$ws1.Range("C13").Value = 2030
$ws1.Range("D13").Value = 189
$ws1.Range("E13").Value = 675
$ws1.Range("F13").Value = 5101
$ws1.Range("B14").Value = 11586
$ws1.Range("C14").Value = 610
$ws1.Range("D14").Value = 2381
$ws1.Range("E14").Value = 3001
$ws1.Range("F14").Value = 5594

# --- Sheet: Noviembre ---
$ws2 = $wb.Worksheets.Item("Noviembre")
$ws2.Range("B2").Value = 2797
$ws2.Range("C2").Value = 1439
$ws2.Range("D2").Value = 200
$ws2.Range("E2").Value = 1158
$ws2.Range("G2").Value = 4.53
$ws2.Range("H2").Value = 12.33
$ws2.Range("B3").Value = 14206
$ws2.Range("C3").Value = 1034
$ws2.Range("D3").Value = 314
$ws2.Range("E3").Value = 959
$ws2.Range("F3").Value = 11899
$ws2.Range("G3").Value = 9.35
$ws2.Range("H3").Value = 22.59
$ws2.Range("B4").Value = 11340
$ws2.Range("C4").Value = 541
$ws2.Range("D4").Value = 266
$ws2.Range("E4").Value = 245
$ws2.Range("F4").Value = 10288
$ws2.Range("G4").Value = 14.14
$ws2.Range("H4").Value = 6.89
$ws2.Range("B5").Value = 20676
$ws2.Range("C5").Value = 2400
$ws2.Range("D5").Value = 65
$ws2.Range("E5").Value = 1286
$ws2.Range("F5").Value = 16925
$ws2.Range("G5").Value = 15.23
$ws2.Range("H5").Value = 4.79
$ws2.Range("I5").Value = 2
$ws2.Range("B6").Value = 9717
$ws2.Range("C6").Value = 2140
$ws2.Range("D6").Value = 162
$ws2.Range("E6").Value = 980
$ws2.Range("F6").Value = 6435
$ws2.Range("G6").Value = 7.97
$ws2.Range("H6").Value = 6.96
$ws2.Range("I6").Value = 9
$ws2.Range("J6").Value = 38
$ws2.Range("B7").Value = 23092
$ws2.Range("C7").Value = 2036
$ws2.Range("D7").Value = 716
$ws2.Range("E7").Value = 2898
$ws2.Range("F7").Value = 17442
$ws2.Range("G7").Value = 26.03
$ws2.Range("H7").Value = 6.33
$ws2.Range("I7").Value = 35
$ws2.Range("J7").Value = 20
$ws2.Range("B8").Value = 12784
$ws2.Range("C8").Value = 921
$ws2.Range("D8").Value = 689
$ws2.Range("E8").Value = 2163
$ws2.Range("F8").Value = 9011
$ws2.Range("G8").Value = 11.89
$ws2.Range("H8").Value = 4.5
$ws2.Range("J8").Value = 96
$ws2.Range("B9").Value = 10570
$ws2.Range("C9").Value = 3992
$ws2.Range("D9").Value = 134
$ws2.Range("E9").Value = 996
$ws2.Range("F9").Value = 5448
$ws2.Range("G9").Value = 10.62
$ws2.Range("H9").Value = 7.72
$ws2.Range("B10").Value = 12697
$ws2.Range("C10").Value = 365
$ws2.Range("D10").Value = 67
$ws2.Range("E10").Value = 274
$ws2.Range("F10").Value = 11991
$ws2.Range("G10").Value = 2.17
$ws2.Range("H10").Value = 17.03
$ws2.Range("I10").Value = 10
$ws2.Range("B11").Value = 15738
$ws2.Range("C11").Value = 1166
$ws2.Range("D11").Value = 557
$ws2.Range("E11").Value = 1768
$ws2.Range("F11").Value = 12247
$ws2.Range("G11").Value = 17.35
$ws2.Range("H11").Value = 5.91
$ws2.Range("I11").Value = 51
$ws2.Range("J11").Value = 116
$ws2.Range("B12").Value = 14363
$ws2.Range("C12").Value = 1635
$ws2.Range("D12").Value = 98
$ws2.Range("E12").Value = 846
$ws2.Range("F12").Value = 11783
$ws2.Range("G12").Value = 8.460000000000001
$ws2.Range("H12").Value = 9.59
$ws2.Range("I12").Value = 21
$ws2.Range("J12").Value = 38
$ws2.Range("B13").Value = 6683
$ws2.Range("C13").Value = 2715
$ws2.Range("D13").Value = 167
$ws2.Range("E13").Value = 1077
$ws2.Range("F13").Value = 2724
$ws2.Range("G13").Value = 7.57
$ws2.Range("H13").Value = 10.03
$ws2.Range("B14").Value = 14206
$ws2.Range("C14").Value = 1034
$ws2.Range("D14").Value = 314
$ws2.Range("E14").Value = 959
$ws2.Range("F14").Value = 11899
$ws2.Range("G14").Value = 9.35
$ws2.Range("H14").Value = 22.59
